# Update "想去人数" (F column) counts for the 展览 (Exhibition) sheet
# and the 全部类型 (All types) sheet, per the upstream data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1): rows 2-9 ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 274
$wsExpo.Range("F3").Value = 165
$wsExpo.Range("F4").Value = 1999
$wsExpo.Range("F5").Value = 1623
$wsExpo.Range("F6").Value = 292
$wsExpo.Range("F7").Value = 79
$wsExpo.Range("F8").Value = 627
$wsExpo.Range("F9").Value = 147

# --- Sheet "全部类型" (sheet4): rows 2-6 and 8-10 (row 7 is a music-concert
#     entry not present in the 展览 sheet, and is left unchanged) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 274
$wsAll.Range("F3").Value = 165
$wsAll.Range("F4").Value = 1999
$wsAll.Range("F5").Value = 1623
$wsAll.Range("F6").Value = 292
$wsAll.Range("F8").Value = 79
$wsAll.Range("F9").Value = 627
$wsAll.Range("F10").Value = 147
